$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows (8-11) below the existing NCAP_START block, each cloned
# from row 7 so that the D:H "empty block" styling (style 10) and the
# J/L border styling come along with the new rows.
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(11).Insert()

# Row-insert drops the hairline border on columns J and L; restore their
# formatting explicitly by re-pasting formats from the still-intact row 7.
$ws.Range("J7").Copy()
$ws.Range("J8:J11").PasteSpecial(-4122)
$ws.Range("L7").Copy()
$ws.Range("L8:L11").PasteSpecial(-4122)

# Materialize the (empty) I and K cells in the new rows so the row matches
# the authored layout, without altering their (default) appearance.
$ws.Range("I8:I11").Interior.Pattern = -4142
$ws.Range("K8:K11").Interior.Pattern = -4142

# Row 8: P-TH-CCGT-GAS-CCS04-Cork1, NCAP_ILED = 0
$ws.Range("D8").Value = "NCAP_ILED"
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = "P-TH-CCGT-GAS-CCS04-Cork1"

# Row 9: P-TH-CCGT-GAS-CCS04-Cork2, NCAP_ILED = 0
$ws.Range("D9").Value = "NCAP_ILED"
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = "P-TH-CCGT-GAS-CCS04-Cork2"

# Row 10: P-TH-CCGT-GAS-CCS04-Dublin1, NCAP_ILED = 0
$ws.Range("D10").Value = "NCAP_ILED"
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = "P-TH-CCGT-GAS-CCS04-Dublin1"

# Row 11: P-TH-CCGT-GAS-CCS04-Dublin2, NCAP_ILED = 0
$ws.Range("D11").Value = "NCAP_ILED"
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = "P-TH-CCGT-GAS-CCS04-Dublin2"

# Match the saved selection state from the authored workbook.
$ws.Range("M16").Select() | Out-Null
